$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 11 new rows right after the existing data row (row 16),
#    pushing the old blank rows + footer block down by 11 rows.
# ------------------------------------------------------------------
for ($i = 0; $i -lt 11; $i++) {
    $ws.Rows.Item(17).Insert(-4121)
}

# ------------------------------------------------------------------
# 2. Copy the formatting (borders/fonts/number formats) of row 16
#    down across the freshly inserted rows 17-27 so the table looks
#    consistent.
# ------------------------------------------------------------------
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Fill in the worker rows. Row 16 now becomes the first worker
#    row, rows 17-26 are the newly added workers, and row 27 gets
#    the original worker (previously row 16) moved to the bottom.
# ------------------------------------------------------------------

function Set-WorkerRow($r, $doc, $name, $periodo, $valorMora, $salario) {
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = $doc
    $ws.Range("D$r").Value = $name
    $ws.Range("E$r").Value = $periodo
    $ws.Range("F$r").Value = $valorMora
    $ws.Range("G$r").Value = $salario
}

Set-WorkerRow 16 "73136979"    "ALFONSO AMADOR SAYAS"              "2507" 1898  1423500
Set-WorkerRow 17 "1128047141"  "RAUL ANTONIO CASTRO JIMENEZ"       "2507" 52000 1300000
Set-WorkerRow 18 "1128047141"  "RAUL ANTONIO CASTRO JIMENEZ"       "2506" 52000 1300000
Set-WorkerRow 19 "8854222"     "DANIEL PATERNINA MENDOZA"          "2507" 1898  1423500
Set-WorkerRow 20 "1087210087"  "JOSE ALEJANDRO SOLARTE RDELO"      "2507" 56940 1423500
Set-WorkerRow 21 "1047509689"  "YERSON ENRIQUE TIJERA RAMIREZ"     "2507" 52000 1300000
Set-WorkerRow 22 "1047509689"  "YERSON ENRIQUE TIJERA RAMIREZ"     "2506" 52000 1300000
Set-WorkerRow 23 "1007469987"  "SEBASTIAN VILLA ORTEGA"            "2506" 56940 1423500
Set-WorkerRow 24 "1151473213"  "FRANCISCO JAVIER VALENCIA COTES"   "2507" 56940 1423500
Set-WorkerRow 25 "1002073302"  "LEIDER ENRIQUE BARRERA NAVARRO"    "2507" 52000 1300000
Set-WorkerRow 26 "1002073302"  "LEIDER ENRIQUE BARRERA NAVARRO"    "2506" 52000 1300000
Set-WorkerRow 27 "1064978493"  "MAURICIO JOSE KARDUSS GONZALEZ"    "2504" 68000 3000000

# ------------------------------------------------------------------
# 4. Update summary header fields: total "Valor Mora" owed, count of
#    workers and count of periods.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 554616
$ws.Range("C13").Value = 9
$ws.Range("F13").Value = 3

# ------------------------------------------------------------------
# 5. Adjust the "best fit" column widths now that longer names and
#    numbers have been added to the table.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.54296875
$ws.Columns.Item(3).ColumnWidth = 16.7265625
$ws.Columns.Item(5).ColumnWidth = 13.54296875
$ws.Columns.Item(6).ColumnWidth = 10.1796875
$ws.Columns.Item(7).ColumnWidth = 14.36328125
$ws.Columns.Item(8).ColumnWidth = 19.36328125
$ws.Columns.Item(9).ColumnWidth = 18.08984375
$ws.Columns.Item(10).ColumnWidth = 15
